$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

# Row 2
Set-TextValue $ws "D2" '73.191.08'
Set-TextValue $ws "E2" '  -0.07%  '

# Row 3
Set-TextValue $ws "D3" '3.977.97'
Set-TextValue $ws "E3" '  -1.69%  '

# Row 4
Set-TextValue $ws "E4" '  -0.07%  '

# Row 5
Set-TextValue $ws "D5" '610.11'

# Row 6
Set-TextValue $ws "D6" '172.31'
Set-TextValue $ws "E6" '  +13.66%  '

# Row 7
Set-TextValue $ws "E7" '  +0.56%  '

# Row 8
Set-TextValue $ws "D8" '1.00'
Set-TextValue $ws "E8" '  +0.02%  '

# Row 9
Set-TextValue $ws "D9" '0.805'
Set-TextValue $ws "E9" '  +6.22%  '

# Row 10
Set-TextValue $ws "E10" '  +9.92%  '

# Row 11
Set-TextValue $ws "D11" '56.67'
Set-TextValue $ws "E11" '  +6.13%  '

# Row 12
Set-TextValue $ws "D12" '0.0000338'
Set-TextValue $ws "E12" '  +4.17%  '

# Row 13
Set-TextValue $ws "D13" '11.77'
Set-TextValue $ws "E13" '  +6.28%  '

# Row 14
Set-TextValue $ws "D14" '4.611.74'
Set-TextValue $ws "E14" '  -1.69%  '

# Row 15
Set-TextValue $ws "D15" '3.983.77'
Set-TextValue $ws "E15" '  -1.67%  '

# Row 16
Set-TextValue $ws "D16" '21.54'
Set-TextValue $ws "E16" '  +3.91%  '

# Row 17
Set-TextValue $ws "E17" '  +0.18%  '

# Row 18
Set-TextValue $ws "E18" '  +1.49%  '

# Row 19
Set-TextValue $ws "D19" '73.068.48'
Set-TextValue $ws "E19" '  -0.15%  '

# Row 20
Set-TextValue $ws "E20" '  -0.70%  '

# Row 21
Set-TextValue $ws "D21" '458.73'
Set-TextValue $ws "E21" '  +4.04%  '

# Row 22
Set-TextValue $ws "E22" '  +5.18%  '

# Row 23
Set-TextValue $ws "E23" '  -0.39%  '

# Row 24
Set-TextValue $ws "E24" '  -4.02%  '

# Row 25
Set-TextValue $ws "D25" '14.33'
Set-TextValue $ws "E25" '  -1.65%  '

# Row 26
Set-TextValue $ws "D26" '4.26'
Set-TextValue $ws "E26" '  -1.35%  '

# Row 27
Set-TextValue $ws "D27" '11.29'
Set-TextValue $ws "E27" '  -2.17%  '

# Row 28
Set-TextValue $ws "E28" '  -2.18%  '

# Row 29
Set-TextValue $ws "D29" '5.88'
Set-TextValue $ws "E29" '  -1.44%  '

# Row 30
Set-TextValue $ws "D30" '36.37'
Set-TextValue $ws "E30" '  -1.57%  '

# Row 31
Set-TextValue $ws "D31" '7.93'
Set-TextValue $ws "E31" '  -0.37%  '

# Row 32
Set-TextValue $ws "D32" '14.10'
Set-TextValue $ws "E32" '  +3.37%  '

# Row 33
Set-TextValue $ws "D33" '49.31'
Set-TextValue $ws "E33" '  +1.69%  '

# Row 34
Set-TextValue $ws "B34" 'PEPE'
Set-TextValue $ws "C34" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws "D34" '0.0000104'
Set-TextValue $ws "E34" '  +17.01%  '

# Row 35
Set-TextValue $ws "B35" 'Hedera'
Set-TextValue $ws "C35" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D35" '0.129'
Set-TextValue $ws "E35" '  -3.20%  '

# Row 36
Set-TextValue $ws "D36" '69.51'
Set-TextValue $ws "E36" '  +3.23%  '

# Row 37
Set-TextValue $ws "D37" '633.05'
Set-TextValue $ws "E37" '  -8.45%  '

# Row 38
Set-TextValue $ws "D38" '0.431'
Set-TextValue $ws "E38" '  -3.42%  '

# Row 39
Set-TextValue $ws "B39" 'ThetaToken'
Set-TextValue $ws "C39" 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws "D39" '3.41'
Set-TextValue $ws "E39" '  +1.12%  '

# Row 40
Set-TextValue $ws "B40" 'Kaspa'
Set-TextValue $ws "C40" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D40" '0.148'
Set-TextValue $ws "E40" '  -1.06%  '

# Row 41
Set-TextValue $ws "E41" '  -0.07%  '

# Row 42
Set-TextValue $ws "D42" '3.35'
Set-TextValue $ws "E42" '  +51.27%  '

# Row 43
Set-TextValue $ws "E43" '  -0.01%  '

# Row 44
Set-TextValue $ws "D44" '0.0488'
Set-TextValue $ws "E44" '  -1.49%  '

# Row 45
Set-TextValue $ws "D45" '10.58'
Set-TextValue $ws "E45" '  -6.15%  '

# Row 46
Set-TextValue $ws "E46" '  -0.26%  '

# Row 47
Set-TextValue $ws "D47" '0.000301'
Set-TextValue $ws "E47" '  +11.01%  '

# Row 48
Set-TextValue $ws "D48" '2.97'
Set-TextValue $ws "E48" '  -11.09%  '

# Row 49
Set-TextValue $ws "E49" '  +1.46%  '

# Row 50
Set-TextValue $ws "E50" '  -4.70%  '

# Row 51
Set-TextValue $ws "B51" 'Maker'
Set-TextValue $ws "C51" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D51" '2.822.88'
Set-TextValue $ws "E51" '  +0.69%  '
